# Fix country-name typos on the "Quick Notes " sheet (column A of the
# country rows) and update the saved view state (scroll position / active
# cell) to match what was left selected when the file was re-uploaded.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quick Notes ")
$ws.Activate()

# Correct the country labels in column A.
$ws.Range("A8").Value  = "South Korea"   # was "South korea"
$ws.Range("A11").Value = "Hong Kong"     # was "Hongkong"
$ws.Range("A14").Value = "India"         # was "India " (trailing space)
$ws.Range("A16").Value = "Vietnam"       # was "Vietnam " (trailing space)

# Update the view state: scrolled so row 11 is at the top, with A12 selected.
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A12").Select()
